$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# these are written with a temporary Text number format so they stay inline strings,
# then the format is reset back to Normal so no visible style/format change remains.
$textForced = @(
    @('D5', '248.19'),
    @('D6', '0.626'),
    @('D7', '59.55'),
    @('D10', '0.0809'),
    @('D12', '15.28'),
    @('D14', '0.855'),
    @('D15', '22.08'),
    @('D16', '5.47'),
    @('D19', '70.55'),
    @('D21', '5.28'),
    @('D22', '230.67'),
    @('D25', '2.36'),
    @('D26', '9.41'),
    @('D27', '164.69'),
    @('D28', '0.138'),
    @('D29', '19.96'),
    @('D32', '0.0683'),
    @('D33', '4.79'),
    @('D35', '4.51'),
    @('D36', '3.56'),
    @('D38', '1.81'),
    @('D39', '5.49'),
    @('D40', '2.99'),
    @('D41', '0.0981'),
    @('D43', '1.19'),
    @('D44', '16.73'),
    @('D45', '92.08'),
    @('D46', '1.07'),
    @('D48', '7.55'),
    @('D49', '2.12')
)

foreach ($entry in $textForced) {
    $cell = $ws.Range($entry[0])
    $cell.NumberFormat = "@"
    $cell.Value = $entry[1]
    $cell.Style = "Normal"
}

# Remaining cells: new text is not a valid numeric literal (thousands-dotted
# prices, subscript-notation prices, URLs, coin names) so a plain .Value assignment
# already round-trips as text with no extra formatting needed.
$ws.Range('D2').Value = '37.415.29'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '2.043.57'
$ws.Range('E3').Value = '  +4.30%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('E7').Value = '  -1.06%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E9').Value = '  +4.63%  '
$ws.Range('E10').Value = '  +3.07%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('E12').Value = '  +8.27%  '
$ws.Range('D13').Value = '2.341.37'
$ws.Range('E13').Value = '  +4.57%  '
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('E16').Value = '  +3.94%  '
$ws.Range('D17').Value = '2.042.40'
$ws.Range('E17').Value = '  +4.51%  '
$ws.Range('D18').Value = '37.377.15'
$ws.Range('E18').Value = '  +2.57%  '
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').Value = '0.0₃0863'
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('E21').Value = '  +4.11%  '
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +5.46%  '
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('E26').Value = '  +3.38%  '
$ws.Range('E27').Value = '  +2.21%  '
$ws.Range('E28').Value = '  -4.32%  '
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('E30').Value = '  +4.93%  '
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('E32').Value = '  +11.62%  '
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  +13.20%  '
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('E36').Value = '  +5.07%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  +1.93%  '
$ws.Range('E39').Value = '  +1.91%  '
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('E42').Value = '  +4.25%  '
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('E44').Value = '  +6.04%  '
$ws.Range('E45').Value = '  +4.19%  '
$ws.Range('E46').Value = '  +4.43%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '1.385.56'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('E48').Value = '  +5.01%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E49').Value = '  +14.91%  '
$ws.Range('E50').Value = '  +2.11%  '
$ws.Range('D51').Value = '2.232.99'
$ws.Range('E51').Value = '  +4.59%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
